# "add pillars scene kpis"
# Rework the Programs sheet:
#  - rename a few header / KPI field names
#  - collapse the per-program columns so brand_name / display_brand_name
#    mirror program_name
#  - rename program values (Coke with Meal -> Coke with Meals,
#    March Madness -> NCAA, Hydration -> Coke Hydration)
#  - drop the "Big Games" and "Spring/ Summer" rows
#  - bold/underline/blacken the two new header cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cosmetic: nudge the sheet-tab-bar split ratio (best effort; some hosts
# don't expose this all the way through to the saved bookView).
try { $wb.Windows.Item(1).TabRatio = 991 } catch { }

# ---- Remove the two trailing KPI rows (Big Games, Spring/ Summer) ----
$ws.Rows("5:6").Delete()

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "program_name"
$ws.Range("B1").Value = "brand_name"
$ws.Range("C1").Value = "display_brand_name"
$ws.Range("D1").Value = "quri_scene_type"
$ws.Range("E1").Value = "survey_question"
$ws.Range("F1").Value = "survey_target_answer"
$ws.Range("G1").Value = "start_date"
$ws.Range("H1").Value = "end_date"

# New distinct look for the two renamed header cells (bold, underlined, black)
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").Font.Underline = $true
$ws.Range("C1:D1").Font.Color = 0

# ---- Data rows ----
# Row 2: Coke with Meals
$ws.Range("A2").Value = "Coke with Meals"
$ws.Range("B2").Value = "Coke with Meals"
$ws.Range("C2").Value = "Coke with Meals"
$ws.Range("D2").Value = "Warm Display with Coke With Meals POS"

# Row 3: NCAA
$ws.Range("A3").Value = "NCAA"
$ws.Range("B3").Value = "NCAA"
$ws.Range("C3").Value = "NCAA"
$ws.Range("D3").Value = "Warm Display with NCAA March Madness POS"

# Row 4: Coke Hydration
$ws.Range("A4").Value = "Coke Hydration"
$ws.Range("B4").Value = "Coke Hydration"
$ws.Range("C4").Value = "Coke Hydration"
$ws.Range("D4").Value = "Warm Display with Hydration POS"

# ---- Row heights (slightly taller to fit the extra column content) ----
$ws.Rows(2).RowHeight = 23.95
$ws.Rows(3).RowHeight = 23.95
$ws.Rows(4).RowHeight = 23.95

# ---- Column widths (tightened up slightly) ----
$ws.Columns(1).ColumnWidth = 19.146258503401366
$ws.Columns(2).ColumnWidth = 25.217687074829968
$ws.Columns(3).ColumnWidth = 22.92687074829937
$ws.Columns(4).ColumnWidth = 22.92687074829937
$ws.Columns(5).ColumnWidth = 13.207482993197265
$ws.Columns(6).ColumnWidth = 17.39115646258507
$ws.Columns(7).ColumnWidth = 7.671768707482998
$ws.Columns(8).ColumnWidth = 7.130952380952377
